$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A26").Value = "Golang developer"
$ws.Range("B26").Value = "https://www.dice.com/job-detail/ccf40fb8-8acf-4fa1-a932-6b7cdd564006"
$ws.Range("C26").Value = "Richmond, Virginia"
$ws.Range("D26").Value = "Contract"
$ws.Range("E26").Value = "Depends on Experience"
$ws.Range("F26").Value = "Lorvenk Technologies LLC"
